# Generate Report for Handoff
# Adds a new tracked file (de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.md) as row 9
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$fileId = "de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.md"
$commitSha = "8d11abe4ca16599bc5519d8c8e7c5e64db9a55c8"
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/" + $commitSha + "/e2e/" + $fileId

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A9").Value = $fileId
$wsOverview.Range("B9").Value = "e2e\" + $fileId
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = ""
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-08-22 02:54:00"
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Range("B9").Font.Underline = 2
$wsOverview.Range("B9").Font.Color = 15570276

$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), $baseUrl, "", "", "e2e\" + $fileId) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A9").Value = $fileId
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "e2e"
$wsZhCn.Range("E9").Value = "ht"
$wsZhCn.Range("F9").Value = "False"
$wsZhCn.Range("G9").Value = "de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.63d3b4067c4d91fb4701731b8107fe798a18b2a4.zh-cn.xlf"
$wsZhCn.Range("H9").Value = "2016-08-22 02:53:56"
$wsZhCn.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I9").Value = ""
$wsZhCn.Range("J9").Value = ""
$wsZhCn.Range("K9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L9").Value = ""
$wsZhCn.Range("M9").Value = "True"
$wsZhCn.Range("N9").Value = ""
$wsZhCn.Range("O9").Value = "False"
$wsZhCn.Range("P9").Value = ""

$wsZhCn.Range("A9").Font.Underline = 2
$wsZhCn.Range("A9").Font.Color = 15570276

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A9"), $baseUrl, "", "", $fileId) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A9").Value = $fileId
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "e2e"
$wsDeDe.Range("E9").Value = "ht"
$wsDeDe.Range("F9").Value = "False"
$wsDeDe.Range("G9").Value = "de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.63d3b4067c4d91fb4701731b8107fe798a18b2a4.de-de.xlf"
$wsDeDe.Range("H9").Value = "2016-08-22 02:54:00"
$wsDeDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I9").Value = ""
$wsDeDe.Range("J9").Value = ""
$wsDeDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L9").Value = ""
$wsDeDe.Range("M9").Value = "True"
$wsDeDe.Range("N9").Value = ""
$wsDeDe.Range("O9").Value = "False"
$wsDeDe.Range("P9").Value = ""

$wsDeDe.Range("A9").Font.Underline = 2
$wsDeDe.Range("A9").Font.Color = 15570276

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A9"), $baseUrl, "", "", $fileId) | Out-Null
